$wb = $excel.ActiveWorkbook

# Map of worksheet name -> F column updates (row => new value)
$updates = @{
    3  = 126
    4  = 67
    5  = 527
    6  = 7232
    7  = 208
    9  = 1055
    10 = 449
    11 = 11
    12 = 152
    13 = 193
    14 = 636
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
